{"js": "// Update the division-facts table: each data cell holds a short\n// \"AA\u00f7B=\" expression. The commit swaps in a freshly generated set of\n// problems (same table shape/formatting, only the operands differ).\n// We target cells by (row, col) and verify the old text before writing\n// the new one, so the script is safe even if some values collide\n// between the old and new sets (e.g. \"59\u00f75=\" is an old value in one\n// cell and a new value in another).\n\nconst replacements = [\n  [0, 0, \"53\u00f78=\", \"59\u00f75=\"],\n  [0, 1, \"74\u00f76=\", \"47\u00f74=\"],\n  [0, 2, \"96\u00f79=\", \"32\u00f78=\"],\n  [0, 3, \"36\u00f75=\", \"97\u00f78=\"],\n  [0, 4, \"87\u00f72=\", \"20\u00f77=\"],\n  [4, 0, \"90\u00f77=\", \"19\u00f76=\"],\n  [4, 1, \"86\u00f75=\", \"89\u00f74=\"],\n  [4, 2, \"59\u00f75=\", \"50\u00f75=\"],\n  [4, 3, \"43\u00f75=\", \"57\u00f79=\"],\n  [4, 4, \"36\u00f73=\", \"51\u00f77=\"],\n  [8, 0, \"49\u00f79=\", \"88\u00f78=\"],\n  [8, 1, \"92\u00f76=\", \"83\u00f79=\"],\n  [8, 2, \"65\u00f75=\", \"12\u00f78=\"],\n  [8, 3, \"75\u00f75=\", \"76\u00f73=\"],\n  [8, 4, \"89\u00f77=\", \"55\u00f78=\"],\n  [12, 0, \"30\u00f76=\", \"71\u00f78=\"],\n  [12, 1, \"73\u00f79=\", \"57\u00f76=\"],\n  [12, 2, \"95\u00f73=\", \"73\u00f74=\"],\n  [12, 3, \"15\u00f76=\", \"71\u00f75=\"],\n  [12, 4, \"55\u00f74=\", \"84\u00f74=\"],\n  [16, 0, \"81\u00f75=\", \"23\u00f79=\"],\n  [16, 1, \"69\u00f78=\", \"45\u00f75=\"],\n  [16, 2, \"79\u00f77=\", \"66\u00f75=\"],\n  [16, 3, \"23\u00f75=\", \"76\u00f73=\"],\n  [16, 4, \"78\u00f73=\", \"15\u00f79=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst cells = replacements.map(([row, col]) => table.getCell(row, col));\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\ncells.forEach((cell, i) => {\n  const [, , oldText, newText] = replacements[i];\n  const current = cell.value;\n  if (current === oldText) {\n    cell.value = newText;\n  } else if (current !== newText) {\n    // Unexpected existing content \u2014 report loudly instead of silently\n    // overwriting the wrong cell.\n    throw new Error(\n      `Cell (${replacements[i][0]}, ${replacements[i][1]}) had ` +\n        `\"${current}\", expected \"${oldText}\"`\n    );\n  }\n});\n\nawait context.sync();\n", "ps1": "# Update the division-facts table: each data cell holds a short\n# \"AA\u00f7B=\" expression. The commit swaps in a freshly generated set of\n# problems (same table shape/formatting, only the operands differ).\n# Cells are addressed by their 1-based (Row, Col) table coordinates and\n# the existing text is verified before overwriting, so the script is\n# safe even though some values collide between the old and new sets\n# (e.g. \"59\u00f75=\" is an old value in one cell and a new value in\n# another).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    [pscustomobject]@{ Row = 1; Col = 1; OldText = \"53\u00f78=\"; NewText = \"59\u00f75=\" },\n    [pscustomobject]@{ Row = 1; Col = 2; OldText = \"74\u00f76=\"; NewText = \"47\u00f74=\" },\n    [pscustomobject]@{ Row = 1; Col = 3; OldText = \"96\u00f79=\"; NewText = \"32\u00f78=\" },\n    [pscustomobject]@{ Row = 1; Col = 4; OldText = \"36\u00f75=\"; NewText = \"97\u00f78=\" },\n    [pscustomobject]@{ Row = 1; Col = 5; OldText = \"87\u00f72=\"; NewText = \"20\u00f77=\" },\n    [pscustomobject]@{ Row = 5; Col = 1; OldText = \"90\u00f77=\"; NewText = \"19\u00f76=\" },\n    [pscustomobject]@{ Row = 5; Col = 2; OldText = \"86\u00f75=\"; NewText = \"89\u00f74=\" },\n    [pscustomobject]@{ Row = 5; Col = 3; OldText = \"59\u00f75=\"; NewText = \"50\u00f75=\" },\n    [pscustomobject]@{ Row = 5; Col = 4; OldText = \"43\u00f75=\"; NewText = \"57\u00f79=\" },\n    [pscustomobject]@{ Row = 5; Col = 5; OldText = \"36\u00f73=\"; NewText = \"51\u00f77=\" },\n    [pscustomobject]@{ Row = 9; Col = 1; OldText = \"49\u00f79=\"; NewText = \"88\u00f78=\" },\n    [pscustomobject]@{ Row = 9; Col = 2; OldText = \"92\u00f76=\"; NewText = \"83\u00f79=\" },\n    [pscustomobject]@{ Row = 9; Col = 3; OldText = \"65\u00f75=\"; NewText = \"12\u00f78=\" },\n    [pscustomobject]@{ Row = 9; Col = 4; OldText = \"75\u00f75=\"; NewText = \"76\u00f73=\" },\n    [pscustomobject]@{ Row = 9; Col = 5; OldText = \"89\u00f77=\"; NewText = \"55\u00f78=\" },\n    [pscustomobject]@{ Row = 13; Col = 1; OldText = \"30\u00f76=\"; NewText = \"71\u00f78=\" },\n    [pscustomobject]@{ Row = 13; Col = 2; OldText = \"73\u00f79=\"; NewText = \"57\u00f76=\" },\n    [pscustomobject]@{ Row = 13; Col = 3; OldText = \"95\u00f73=\"; NewText = \"73\u00f74=\" },\n    [pscustomobject]@{ Row = 13; Col = 4; OldText = \"15\u00f76=\"; NewText = \"71\u00f75=\" },\n    [pscustomobject]@{ Row = 13; Col = 5; OldText = \"55\u00f74=\"; NewText = \"84\u00f74=\" },\n    [pscustomobject]@{ Row = 17; Col = 1; OldText = \"81\u00f75=\"; NewText = \"23\u00f79=\" },\n    [pscustomobject]@{ Row = 17; Col = 2; OldText = \"69\u00f78=\"; NewText = \"45\u00f75=\" },\n    [pscustomobject]@{ Row = 17; Col = 3; OldText = \"79\u00f77=\"; NewText = \"66\u00f75=\" },\n    [pscustomobject]@{ Row = 17; Col = 4; OldText = \"23\u00f75=\"; NewText = \"76\u00f73=\" },\n    [pscustomobject]@{ Row = 17; Col = 5; OldText = \"78\u00f73=\"; NewText = \"15\u00f79=\" }\n)\n\nforeach ($r in $replacements) {\n    $cell = $t.Cell($r.Row, $r.Col)\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -eq $r.OldText) {\n        $cell.Range.Text = $r.NewText\n    } elseif ($current -ne $r.NewText) {\n        throw \"Cell ($($r.Row), $($r.Col)) had '$current', expected '$($r.OldText)'\"\n    }\n}\n"}
